$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Encl. door sensor" entry: rows 15-16 in the Optional sensors table ---
# Description (merged F15:F16)
$ws.Range("F15").Value = "Encl. door sensor"
$ws.Range("F15").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F15").VerticalAlignment = -4108     # xlCenter
$ws.Range("F16").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F16").VerticalAlignment = -4108     # xlCenter

# Channel column
$ws.Range("G15").Value = "state detect"
$ws.Range("G16").Value = "excitation"

# Color column (reuses the existing "white" shared string)
$ws.Range("H15").Value = "white"
$ws.Range("H16").Value = "white"

# CR3000 column
$ws.Range("I15").Value = "C8"
$ws.Range("I16").Value = "5V"

# Merge the description cell across the two new rows
$ws.Range("F15:F16").Merge()

# --- Update the saved selection to match the authored edit ---
$ws.Range("F17").Select()
